# "Error Calculations and Plots"
# The sheet holds a missing-data table (ID + columns A-D,F). This edit:
#   1. Removes two rows entirely ("RM 232" and "SC 92"), shifting all
#      subsequent rows up (35 data rows -> 33 data rows).
#   2. Re-marks a scattered set of cells as missing (blank) or fills in
#      previously-missing cells with newly imputed numeric values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Delete the two removed rows -----------------------------------
# Row 26 ("RM 232") is deleted first; afterwards the old row 28
# ("SC 92") has shifted up to row 27, so it is deleted there.
$ws.Rows.Item(26).Delete()
$ws.Rows.Item(27).Delete()

# --- 2. Apply the scattered value / blank-out changes ------------------
$ws.Range("E2").Value = ""
$ws.Range("F4").Value = ""
$ws.Range("E5").Value = -5
$ws.Range("D6").Value = -14.2
$ws.Range("E6").Value = -5.7
$ws.Range("F6").Value = 16.43
$ws.Range("D8").Value = ""
$ws.Range("E9").Value = ""
$ws.Range("E10").Value = ""
$ws.Range("F11").Value = 17.65
$ws.Range("D12").Value = -14.1
$ws.Range("F12").Value = ""
$ws.Range("D14").Value = ""
$ws.Range("F14").Value = 17.76
$ws.Range("D17").Value = -14.7
$ws.Range("F17").Value = ""
$ws.Range("D18").Value = -15.2
$ws.Range("D19").Value = ""
$ws.Range("F19").Value = 17.81
$ws.Range("D20").Value = ""
$ws.Range("F21").Value = 16.58
$ws.Range("F22").Value = 16.81
$ws.Range("D23").Value = -13.9
$ws.Range("E24").Value = -8.1
$ws.Range("F25").Value = ""
$ws.Range("F26").Value = ""
$ws.Range("B27").Value = -20.4
$ws.Range("D27").Value = ""
$ws.Range("F27").Value = ""
$ws.Range("B28").Value = ""
$ws.Range("E28").Value = ""
$ws.Range("F28").Value = ""
$ws.Range("B29").Value = ""
$ws.Range("B30").Value = -19.7
$ws.Range("E30").Value = -5.7
$ws.Range("F31").Value = 17.18
$ws.Range("B32").Value = ""
